$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: Snapshot the cell-style patterns we will need later onto scratch
# rows far below the data (so later row deletions inside rows 12-28 don't
# disturb them). xlPasteFormats (-4122) copies only formatting, not values.
#   Row 19  -> pattern "label in A + value in B/C"                 (h=120)
#   Row 20  -> pattern "label only in A"                           (h=120)
#   Row 27  -> pattern "value only in B/C"                         (h=30)
# ---------------------------------------------------------------------------
$ws.Range("A19:C19").Copy() | Out-Null
$ws.Range("A500:C500").PasteSpecial(-4122)

$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A501:C501").PasteSpecial(-4122)

$ws.Range("A27:C27").Copy() | Out-Null
$ws.Range("A502:C502").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: Row 10 keeps its "Objetivos:" label; only the long paragraph in
# B10/C10 is replaced by the professor name that used to sit further down.
# ---------------------------------------------------------------------------
$ws.Range("B10:C10").Value = "5840726 - Cristina Bormio Nunes"

# Row 11 ("Objectives:") is untouched.

# ---------------------------------------------------------------------------
# Step 3: Remove the old rows 12-28 entirely; we rebuild rows 12-23 from
# scratch right after. (This also shifts the scratch rows 500-502 up to
# 483-485.)
# ---------------------------------------------------------------------------
$ws.Range("A12:C28").EntireRow.Delete()

$rowAB60  = 10   # template: label A + value B/C, height 60
$rowA60   = 11   # template: label A only, height 60
$rowAB120 = 483  # template: label A + value B/C, height 120 (was row 19)
$rowA120  = 484  # template: label A only, height 120 (was row 20)
$rowBC30  = 485  # template: value B/C only, height 30 (was row 27)

# ---------------------------------------------------------------------------
# Row 12: "Programa resumido:" label + "6495737 - Durval Rodrigues Junior"
# ---------------------------------------------------------------------------
$ws.Range("A$rowAB60`:C$rowAB60").Copy() | Out-Null
$ws.Range("A12:C12").PasteSpecial(-4122)
$ws.Rows(12).RowHeight = 60
$ws.Range("A12").Value = "Programa resumido:"
$ws.Range("B12:C12").Value = "6495737 - Durval Rodrigues Junior"

# ---------------------------------------------------------------------------
# Row 13: "Short syllabus:" label only
# ---------------------------------------------------------------------------
$ws.Range("A$rowA60`:C$rowA60").Copy() | Out-Null
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Rows(13).RowHeight = 60
$ws.Range("A13").Value = "Short syllabus:"

# ---------------------------------------------------------------------------
# Row 14: "Programa:" label + "1341653 - Maria José Ramos Sandim"
# ---------------------------------------------------------------------------
$ws.Range("A$rowAB120`:C$rowAB120").Copy() | Out-Null
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Rows(14).RowHeight = 120
$ws.Range("A14").Value = "Programa:"
$ws.Range("B14:C14").Value = "1341653 - Maria José Ramos Sandim"

# ---------------------------------------------------------------------------
# Row 15: "Syllabus:" label only
# ---------------------------------------------------------------------------
$ws.Range("A$rowA120`:C$rowA120").Copy() | Out-Null
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Rows(15).RowHeight = 120
$ws.Range("A15").Value = "Syllabus:"

# ---------------------------------------------------------------------------
# Row 16: "Avaliação:" label only, default height
# ---------------------------------------------------------------------------
$ws.Range("A$rowA120`:C$rowA120").Copy() | Out-Null
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Avaliação:"

# ---------------------------------------------------------------------------
# Row 17: "Método:" label + "1643715 - Paulo Atsushi Suzuki"
# ---------------------------------------------------------------------------
$ws.Range("A$rowAB60`:C$rowAB60").Copy() | Out-Null
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Rows(17).RowHeight = 60
$ws.Range("A17").Value = "Método:"
$ws.Range("B17:C17").Value = "1643715 - Paulo Atsushi Suzuki"

# ---------------------------------------------------------------------------
# Row 18: "Critério:" label + "Aulas expositivas teóricas, aulas de exercícios."
# ---------------------------------------------------------------------------
$ws.Range("A$rowAB60`:C$rowAB60").Copy() | Out-Null
$ws.Range("A18:C18").PasteSpecial(-4122)
$ws.Rows(18).RowHeight = 60
$ws.Range("A18").Value = "Critério:"
$ws.Range("B18:C18").Value = "Aulas expositivas teóricas, aulas de exercícios."

# ---------------------------------------------------------------------------
# Row 19: "Norma de recuperação:" label + "Duas provas escritas: ..."
# ---------------------------------------------------------------------------
$ws.Range("A$rowAB60`:C$rowAB60").Copy() | Out-Null
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Rows(19).RowHeight = 60
$ws.Range("A19").Value = "Norma de recuperação:"
$ws.Range("B19:C19").Value = "Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + 2P2)/3"

# ---------------------------------------------------------------------------
# Row 20: "Bibliografia:" label + the recovery-exam paragraph
# ---------------------------------------------------------------------------
$ws.Range("A$rowAB120`:C$rowAB120").Copy() | Out-Null
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Rows(20).RowHeight = 120
$ws.Range("A20").Value = "Bibliografia:"
$ws.Range("B20:C20").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# ---------------------------------------------------------------------------
# Row 21: "Requisitos:" label only, default height
# ---------------------------------------------------------------------------
$ws.Range("A$rowA120`:C$rowA120").Copy() | Out-Null
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Requisitos:"

# ---------------------------------------------------------------------------
# Row 22/23: requisite lines, value only (B/C), height 30
# ---------------------------------------------------------------------------
$ws.Range("A$rowBC30`:C$rowBC30").Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Rows(22).RowHeight = 30
$ws.Range("B22:C22").Value = "LOB1003 -  Cálculo I  (Requisito)`n"

$ws.Range("A$rowBC30`:C$rowBC30").Copy() | Out-Null
$ws.Range("A23:C23").PasteSpecial(-4122)
$ws.Rows(23).RowHeight = 30
$ws.Range("B23:C23").Value = "LOB1006 -  Cálculo IV  (Requisito)`n"

# ---------------------------------------------------------------------------
# Step 4: Clean up the scratch template rows (483-485) so they don't show up
# in the saved sheet / dimension.
# ---------------------------------------------------------------------------
$ws.Range("A483:C485").ClearContents()
$ws.Range("A483:C485").ClearFormats()

Write-Output ("Final UsedRange: " + $ws.UsedRange.Address())
